$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns (AC, AD, AE), matching the
# style used by the rest of the header row (row 1).
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header
# cell onto the three new header cells.
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null

# Restore the text we just set (PasteSpecial formats only, but make sure
# values remain correct).
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team record (Wins=91, Losses=71, Ties=0) for every data row.
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 91
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 0
}
